$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Rushing" - Week 17 stat updates (no new rows)
# ---------------------------------------------------------------------------
$rush = $wb.Worksheets.Item("Rushing")

$rush.Range("D2").Value = 9                                        # M.Jones

$rush.Range("C3").Value = 117                                      # D.Harris
$rush.Range("D3").Value = 59
$rush.Range("E3").Value = 16
$rush.Range("F3").Value = 43

$rush.Range("C5").Value = 79                                       # R.Stevenson
$rush.Range("D5").Value = 43
$rush.Range("F5").Value = 25

$rush.Range("C7").Value = 7                                        # B.Bolden

$rush.Range("C10").Value = 8                                       # K.Bourne
$rush.Range("D10").Value = 4
$rush.Range("F10").Value = 2

$rush.Range("D13").Value = 5                                       # J.Smith

# ---------------------------------------------------------------------------
# Sheet "Receiving" - Week 17 stat updates + newly logged player K.Wilkerson
# ---------------------------------------------------------------------------
$recv = $wb.Worksheets.Item("Receiving")

$recv.Range("C2").Value = 17                                       # D.Harris
$recv.Range("D2").Value = 13

$recv.Range("C5").Value = 41                                       # B.Bolden
$recv.Range("D5").Value = 33

$recv.Range("C7").Value = 94                                       # J.Meyers
$recv.Range("D7").Value = 70
$recv.Range("G7").Value = 13
$recv.Range("H7").Value = 8

$recv.Range("C8").Value = 51                                       # K.Bourne
$recv.Range("D8").Value = 42
$recv.Range("E8").Value = 14
$recv.Range("F8").Value = 10
$recv.Range("G8").Value = 6
$recv.Range("H8").Value = 4

$recv.Range("C9").Value = 2                                        # G.Olszewski
$recv.Range("D9").Value = 1

# Insert a row for the newly logged K.Wilkerson before J.Smith, pushing
# J.Smith down a row and giving H.Henry a brand new row at the end, so the
# final order becomes: K.Wilkerson, J.Smith, H.Henry.
#
# Row 14 does not exist yet, so clone row 13's "index column" formatting
# into it first (A13/A12 already carry the bordered style and keep it
# automatically when their .Value is simply overwritten below).
$recv.Range("A13").Copy()
$recv.Range("A14").PasteSpecial(-4122)     # xlPasteFormats

# Row 12 becomes the newly logged player K.Wilkerson.
$recv.Range("B12").Value = "K.Wilkerson"
$recv.Range("C12").Value = 5
$recv.Range("D12").Value = 3
$recv.Range("E12").Value = 3
$recv.Range("F12").Value = 1
$recv.Range("G12").Value = 2
$recv.Range("H12").Value = 2

# Row 13 becomes J.Smith (previously row 12), with updated Week 17 stats.
$recv.Range("A13").Value = 11
$recv.Range("B13").Value = "J.Smith"
$recv.Range("C13").Value = 34
$recv.Range("D13").Value = 23
$recv.Range("E13").Value = 11
$recv.Range("F13").Value = 5
$recv.Range("G13").Value = 10
$recv.Range("H13").Value = 5

# Row 14 is a brand-new row for H.Henry (previously row 13).
$recv.Range("A14").Value = 12
$recv.Range("B14").Value = "H.Henry"
$recv.Range("C14").Value = 58
$recv.Range("D14").Value = 40
$recv.Range("E14").Value = 12
$recv.Range("F14").Value = 6
$recv.Range("G14").Value = 18
$recv.Range("H14").Value = 9
